$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Previously added")
$ws2 = $wb.Worksheets.Item("New")

# ---------------------------------------------------------------------------
# 1) "Previously added" sheet: append the entry that used to be the sole row
#    of the "New" sheet (row 2) as the new last row (row 322). It keeps the
#    same link/price/district/area/cadastre/date it already had.
# ---------------------------------------------------------------------------
$ws1.Range("A322").Value = "https://www.ss.com/msg/lv/real-estate/wood/limbadzi-and-reg/ledurgas-pag/obmbf.html"
$ws1.Range("B322").Value = "76 000 €"
$ws1.Range("C322").Value = "Limbaži un raj."
$ws1.Range("D322").Value = "35 ha."
$ws1.Range("E322").Value = "66560010200"
$ws1.Range("F322").Value = 45995.77083333333
$ws1.Hyperlinks.Add($ws1.Range("A322"), "https://www.ss.com/msg/lv/real-estate/wood/limbadzi-and-reg/ledurgas-pag/obmbf.html")

# Adding the hyperlink resets the cell style of A322, so re-apply the usual
# row formatting (hyperlink style in A, plain body style in B:E, date style
# in F) by copying it over from the row right above, which already has it.
$ws1.Range("A321:F321").Copy()
$ws1.Range("A322:F322").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) "New" sheet: the old row 2 (now duplicated onto "Previously added")
#    is replaced by fresh scraped listings, and five more new rows are
#    appended below it.
# ---------------------------------------------------------------------------

# Drop the stale hyperlink that used to sit on A2 before overwriting it.
$ws2.Range("A2").Hyperlinks.Delete()

$ws2.Range("A2").Value = "https://www.ss.com/msg/lv/real-estate/wood/balvi-and-reg/vecumu-pag/cogdb.html"
$ws2.Range("B2").Value = "58 000 €"
$ws2.Range("C2").Value = "Balvi un raj."
$ws2.Range("D2").Value = "19.47 ha."
$ws2.Range("E2").Value = "38920050109"
$ws2.Range("F2").Value = 45999.48402777778
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://www.ss.com/msg/lv/real-estate/wood/balvi-and-reg/vecumu-pag/cogdb.html")

$ws2.Range("A3").Value = "https://www.ss.com/msg/lv/real-estate/wood/balvi-and-reg/vecumu-pag/amgxp.html"
$ws2.Range("B3").Value = "43 000 €"
$ws2.Range("C3").Value = "Balvi un raj."
$ws2.Range("D3").Value = "14.46 ha."
$ws2.Range("E3").Value = "38920050232"
$ws2.Range("F3").Value = 45999.45486111111
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://www.ss.com/msg/lv/real-estate/wood/balvi-and-reg/vecumu-pag/amgxp.html")

$ws2.Range("A4").Value = "https://www.ss.com/msg/lv/real-estate/wood/bauska-and-reg/iecavas-nov/ijmhp.html"
$ws2.Range("B4").Value = "32 100 €"
$ws2.Range("C4").Value = "Bauska un raj."
$ws2.Range("D4").Value = "3 ha."
$ws2.Range("E4").Value = "40640130163"
$ws2.Range("F4").Value = 45998.55972222222
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://www.ss.com/msg/lv/real-estate/wood/bauska-and-reg/iecavas-nov/ijmhp.html")

$ws2.Range("A5").Value = "https://www.ss.com/msg/lv/real-estate/wood/jekabpils-and-reg/aknistes-l-t/hfxcb.html"
$ws2.Range("B5").Value = "60 000 €"
$ws2.Range("C5").Value = "Jēkabpils un raj."
$ws2.Range("D5").Value = "16.63 ha."
$ws2.Range("E5").Value = "56250070275"
$ws2.Range("F5").Value = 45999.459027777775
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://www.ss.com/msg/lv/real-estate/wood/jekabpils-and-reg/aknistes-l-t/hfxcb.html")

$ws2.Range("A6").Value = "https://www.ss.com/msg/lv/real-estate/wood/ludza-and-reg/brigu-pag/kniic.html"
$ws2.Range("B6").Value = "32 000 €"
$ws2.Range("C6").Value = "Ludza un raj."
$ws2.Range("D6").Value = "10.21 ha."
$ws2.Range("E6").Value = "68460030011"
$ws2.Range("F6").Value = 45999.48888888889
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://www.ss.com/msg/lv/real-estate/wood/ludza-and-reg/brigu-pag/kniic.html")

$ws2.Range("A7").Value = "https://www.ss.com/msg/lv/real-estate/wood/ludza-and-reg/merdzenes-pag/hmobk.html"
$ws2.Range("B7").Value = "20 000 €"
$ws2.Range("C7").Value = "Ludza un raj."
$ws2.Range("D7").Value = "5.60 ha."
$ws2.Range("E7").Value = "68720020331"
$ws2.Range("F7").Value = 45999.46319444444
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://www.ss.com/msg/lv/real-estate/wood/ludza-and-reg/merdzenes-pag/hmobk.html")

# Restore the normal row formatting (hyperlink style in A, body style in
# B:E, date style in F) across rows 2-7 in one shot, since Hyperlinks.Add
# stomps on the style of the cells it touches.
$ws1.Range("A321:F321").Copy()
$ws2.Range("A2:F7").PasteSpecial(-4122)

Write-Output "done"
